$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A, rows 2..97 currently hold "q1".."q96".
# Renumber them down by one: "q0".."q95".
for ($r = 2; $r -le 97; $r++) {
    $n = $r - 2
    $ws.Cells.Item($r, 1).Value = "q$n"
}
